# Update stats for 2026-01 (row 26) to reflect latest figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B26").Value = 6494
$ws.Range("C26").Value = 1010
$ws.Range("D26").Value = 6049073
$ws.Range("E26").Value = 931.4864490298737
$ws.Range("F26").Value = 9.751563292208898
$ws.Range("G26").Value = 7.218683651804669
$ws.Range("H26").Value = 25.97214901644407
